$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Informe")

# Insert a new column before column E ("Fin"), shifting Recepcion..Recorrido right by one.
$ws.Columns("E:E").Insert()

# New header cell for the inserted column.
$ws.Range("E13").Value = "Fin"

# Update the cell selection to match the saved view state.
$ws.Range("F9").Select()
